$wb = $excel.ActiveWorkbook

# Sheet: y_fitted_on_begin_2016 (sheet1) - fix B2:B24 values
$ws1 = $wb.Worksheets.Item("y_fitted_on_begin_2016")
$ws1.Range("B2").Value = -0.1572072984793992
$ws1.Range("B3").Value = 23.31200298421016
$ws1.Range("B4").Value = 23.2568557127425
$ws1.Range("B5").Value = 23.16062652841528
$ws1.Range("B6").Value = 23.11260328765611
$ws1.Range("B7").Value = 22.79717727268367
$ws1.Range("B8").Value = 22.41344572282036
$ws1.Range("B9").Value = 22.2784483867644
$ws1.Range("B10").Value = 21.45733869698017
$ws1.Range("B11").Value = 20.87463444728024
$ws1.Range("B12").Value = 20.74177588995353
$ws1.Range("B13").Value = 20.80161235211317
$ws1.Range("B14").Value = 21.00064494873929
$ws1.Range("B15").Value = 21.07184045188594
$ws1.Range("B16").Value = 21.37509646969307
$ws1.Range("B17").Value = 21.1599834330719
$ws1.Range("B18").Value = 19.79252546058449
$ws1.Range("B19").Value = 20.19318939152069
$ws1.Range("B20").Value = 19.65918744485145
$ws1.Range("B21").Value = 19.5062893231311
$ws1.Range("B22").Value = 19.24752153567477
$ws1.Range("B23").Value = 19.28860656530033
$ws1.Range("B24").Value = 19.80396922347703

# Sheet: y_fitted_on_begin_2021 (sheet3) - fix B2:B29 values
$ws3 = $wb.Worksheets.Item("y_fitted_on_begin_2021")
$ws3.Range("B2").Value = -0.1557464113639517
$ws3.Range("B3").Value = 23.3134652677244
$ws3.Range("B4").Value = 23.2664529400967
$ws3.Range("B5").Value = 23.17253347678775
$ws3.Range("B6").Value = 23.13222517795666
$ws3.Range("B7").Value = 22.78488770183073
$ws3.Range("B8").Value = 22.38089578835293
$ws3.Range("B9").Value = 22.27136822539923
$ws3.Range("B10").Value = 21.36922367061542
$ws3.Range("B11").Value = 20.78927443152552
$ws3.Range("B12").Value = 20.71639348894494
$ws3.Range("B13").Value = 20.82261216670749
$ws3.Range("B14").Value = 21.05622520113417
$ws3.Range("B15").Value = 21.12282014279839
$ws3.Range("B16").Value = 21.45489781281199
$ws3.Range("B17").Value = 21.181929731297
$ws3.Range("B18").Value = 19.64295437152269
$ws3.Range("B19").Value = 20.21519665085629
$ws3.Range("B20").Value = 19.61896640140471
$ws3.Range("B21").Value = 19.49416213308774
$ws3.Range("B22").Value = 19.23139280369853
$ws3.Range("B23").Value = 19.31040292891043
$ws3.Range("B24").Value = 19.80396922347703
$ws3.Range("B25").Value = 19.99621177421012
$ws3.Range("B26").Value = 20.12546209820798
$ws3.Range("B27").Value = 20.29008194731416
$ws3.Range("B28").Value = 19.98442369411245
$ws3.Range("B29").Value = 19.67530049335799
